# The run's test milestones ("test account milestone" / "test memo milestone")
# need a second, distinct instance so the traceback panel that now renders
# front-and-center on a new run has two rows to show per milestone sheet.
#
# AccountMilestones!A3 and MemoMilestones!A3 currently duplicate row 2's
# milestone name ("test account milestone" / "test memo milestone"); rename
# them to the "...2" variant that CompositeMilestones already references.

$wb = $excel.ActiveWorkbook

$wsAccount = $wb.Worksheets.Item("AccountMilestones")
$wsAccount.Range("A3").Value = "test account milestone 2"

$wsMemo = $wb.Worksheets.Item("MemoMilestones")
$wsMemo.Range("A3").Value = "test memo milestone 2"
